$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
Write-Output "before: $($d.Name)"
$d.Name = "Office Theme"
Write-Output "after (same call): $($d.Name)"
